{"js": "// Swap the body text of the \"Objetivos\" paragraphs with the body text of\n// the \"Programa resumido\" paragraphs (PT line <-> PT line, EN/italic line\n// <-> EN/italic line), keeping each paragraph's own formatting in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two section headings by their exact text.\nlet objetivosIdx = -1;\nlet programaResumidoIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (t === \"Objetivos\") {\n    objetivosIdx = i;\n  } else if (t === \"Programa resumido\") {\n    programaResumidoIdx = i;\n  }\n}\n\nif (objetivosIdx === -1 || programaResumidoIdx === -1) {\n  throw new Error(\"Could not locate 'Objetivos' / 'Programa resumido' headings.\");\n}\n\n// The two paragraphs right after each heading hold the PT text (plain) and\n// the EN text (italic), respectively.\nconst objPt = items[objetivosIdx + 1];\nconst objEn = items[objetivosIdx + 2];\nconst progPt = items[programaResumidoIdx + 1];\nconst progEn = items[programaResumidoIdx + 2];\n\nconst objPtText = objPt.text;\nconst objEnText = objEn.text;\nconst progPtText = progPt.text;\nconst progEnText = progEn.text;\n\n// Swap the text content between the two pairs of paragraphs; formatting\n// (e.g. the italic run properties on the EN paragraphs) stays with the\n// paragraph/run it was already on.\nobjPt.insertText(progPtText, \"Replace\");\nobjEn.insertText(progEnText, \"Replace\");\nprogPt.insertText(objPtText, \"Replace\");\nprogEn.insertText(objEnText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Swap the body text of the \"Objetivos\" paragraphs with the body text of\n# the \"Programa resumido\" paragraphs (PT line <-> PT line, EN/italic line\n# <-> EN/italic line), keeping each paragraph's own formatting in place.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$objetivosIdx = -1\n$programaResumidoIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Objetivos\") {\n        $objetivosIdx = $i\n    } elseif ($t -eq \"Programa resumido\") {\n        $programaResumidoIdx = $i\n    }\n}\n\nif ($objetivosIdx -eq -1 -or $programaResumidoIdx -eq -1) {\n    throw \"Could not locate 'Objetivos' / 'Programa resumido' headings.\"\n}\n\n# The two paragraphs right after each heading hold the PT text (plain) and\n# the EN text (italic), respectively.\n$objPt = $d.Paragraphs.Item($objetivosIdx + 1).Range\n$objEn = $d.Paragraphs.Item($objetivosIdx + 2).Range\n$progPt = $d.Paragraphs.Item($programaResumidoIdx + 1).Range\n$progEn = $d.Paragraphs.Item($programaResumidoIdx + 2).Range\n\n# `.Range.Text` on a paragraph includes the trailing paragraph-mark\n# character (CR, chr 13); strip it so re-assigning the captured text into\n# another paragraph's range doesn't insert an extra paragraph break.\n$objPtText = $objPt.Text.TrimEnd([char]13)\n$objEnText = $objEn.Text.TrimEnd([char]13)\n$progPtText = $progPt.Text.TrimEnd([char]13)\n$progEnText = $progEn.Text.TrimEnd([char]13)\n\n# Swap the text content between the two pairs of paragraphs; formatting\n# (e.g. the italic run properties on the EN paragraphs) stays with the\n# paragraph/run it was already on.\n$objPt.Text = $progPtText\n$objEn.Text = $progEnText\n$progPt.Text = $objPtText\n$progEn.Text = $objEnText\n"}
